$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 94
$ws.Range("I2").Value = 233
$ws.Range("J2").Value = 979
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 243
$ws.Range("M2").Value = 15
$ws.Range("N2").Value = 159
$ws.Range("O2").Value = 1
$ws.Range("R2").Value = 10
$ws.Range("S2").Value = 106
$ws.Range("T2").Value = 158
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 1482
$ws.Range("X2").Value = 1508
$ws.Range("Z2").Value = 27
$ws.Range("AA2").Value = 11
